$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Resumen")

# Use a temporary scratch column far outside the used range to pre-stage the two
# cell styles we need (percent format w/ quote-prefix, then plain w/ quote-prefix)
# so that when the real cells are written below, they land on the same style
# indexes / ordering that a human editing this in Excel would produce.
$ws.Range("Z1").Value = "'1%"
$ws.Range("Z1").NumberFormat = "0.00%"
$ws.Range("Z2").Value = "'x"

# NPS (row 18) changes from 81.32% to 81.16%, kept as text (quote-prefixed) with
# a percentage number format.
$ws.Range("B18:D18").Value = "'81.16%"
$ws.Range("B18:D18").NumberFormat = "0.00%"

# ICX (row 17) changes from 4.8 to 4.80, kept as text (quote-prefixed).
$ws.Range("B17:D17").Value = "'4.80"

# Remove the scratch column entirely so no trace of it remains.
$ws.Columns("Z").Delete()

# Update the active selection to D17, matching where the edit was made.
$ws.Range("D17").Select()
